# Add data for 2022-02-20 (update "through February 12" counts)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header string to reflect the new "as of" date
$ws.Name = "Through 2022-02-12"
$ws.Range("B1").Value = "February 2022 (through February 12)"

# Austin (row 3): B3 2->3, D3 6->7, new P3=1
$ws.Range("B3").Value = 3
$ws.Range("D3").Value = 7
$ws.Range("P3").Value = 1

# New City (row 4): new F4=1
$ws.Range("F4").Value = 1

# South Shore (row 6): L6 2->3
$ws.Range("L6").Value = 3

# Auburn Gresham (row 7): new B7=1, new H7=1
$ws.Range("B7").Value = 1
$ws.Range("H7").Value = 1

# Chatham (row 23): new D23=1
$ws.Range("D23").Value = 1

# Roseland (row 27): new H27=1
$ws.Range("H27").Value = 1

# Portage Park (row 30): new F30=1
$ws.Range("F30").Value = 1

# West Loop (row 34): new L34=1
$ws.Range("L34").Value = 1

# Little Village (row 47): new N47=1
$ws.Range("N47").Value = 1

# Lake View (row 48): new D48=1
$ws.Range("D48").Value = 1

# Ashburn (row 56): new L56=1
$ws.Range("L56").Value = 1

# Chinatown (row 60): B60 1->2
$ws.Range("B60").Value = 2

# Jefferson Park (row 71): new B71=1
$ws.Range("B71").Value = 1

# Lincoln Park (row 72): new J72=1
$ws.Range("J72").Value = 1

# Washington Park (row 84): new N84=1
$ws.Range("N84").Value = 1
